{"js": "// Minor changes to the report & specs:\n//   1. Title \"Lab 4 requirements and specs\" -> \"Lab 3 requirements and specs\"\n//      (split into two runs: \"Lab 3\" and \" requirements and specs\")\n//   2. The stray \"_GoBack\" bookmark that sat at the end of the \"Task Control\n//      Blocks...\" bullet moves into the \"support\" bullet, splitting\n//      \"The task will have support the following user inputs:\" into\n//      \"The task will have su\" | \"pport the following user inputs:\"\n\nconst body = context.document.body;\n\n// --- 1) Title: \"Lab 4\" -> \"Lab 3\", producing two separate runs ---------\nconst titleResults = body.search(\"Lab 4\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\n\nconst titleHit = titleResults.items[0];\ntitleHit.insertText(\"Lab 3\", \"Replace\");\nawait context.sync();\n\n// Re-search for the replaced text to get a fresh, valid range, then drop a\n// bookmark at its end and immediately remove it again \u2014 inserting /\n// deleting a bookmark is what forces Word to keep the text as two runs\n// instead of re-merging them.\nconst newTitleResults = body.search(\"Lab 3\", { matchCase: true });\nnewTitleResults.load(\"items\");\nawait context.sync();\n\nconst newTitleHit = newTitleResults.items[0];\nconst titleSplitRange = newTitleHit.getRange(\"End\");\ntitleSplitRange.insertBookmark(\"_TitleSplit\");\nawait context.sync();\n\ncontext.document.deleteBookmark(\"_TitleSplit\");\nawait context.sync();\n\n// --- 2) Relocate the \"_GoBack\" bookmark -------------------------------\n// Remove the existing (hidden) \"_GoBack\" bookmark from wherever it is.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Find the \"su\" | \"pport\" split point inside the support bullet and drop\n// the bookmark back in at that position.\nconst supportResults = body.search(\"The task will have su\", { matchCase: true });\nsupportResults.load(\"items\");\nawait context.sync();\n\nconst supportHit = supportResults.items[0];\nconst newBookmarkRange = supportHit.getRange(\"End\");\nnewBookmarkRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Minor changes to the report & specs:\n#   1. Title \"Lab 4 requirements and specs\" -> \"Lab 3 requirements and specs\"\n#      (split into two runs: \"Lab 3\" and \" requirements and specs\")\n#   2. The stray \"_GoBack\" bookmark that sat at the end of the \"Task Control\n#      Blocks...\" bullet moves into the \"support\" bullet, splitting\n#      \"The task will have support the following user inputs:\" into\n#      \"The task will have su\" | \"pport the following user inputs:\"\n\n$d = $word.ActiveDocument\n\n# --- 1) Title: \"Lab 4\" -> \"Lab 3\", producing two separate runs ---------\n$d.Content.Find.Execute(\"Lab 4\", $false, $false, $false, $false, $false, $true, 1, $false, \"Lab 3\", 2) | Out-Null\n\n$titlePara = $d.Paragraphs(1)\n$splitAt = $titlePara.Range.Start + (\"Lab 3\").Length\n$titleSplitRange = $d.Range($splitAt, $splitAt)\n$d.Bookmarks.Add(\"_TitleSplit\", $titleSplitRange) | Out-Null\n$d.Bookmarks(\"_TitleSplit\").Delete()\n\n# --- 2) Relocate the \"_GoBack\" bookmark -------------------------------\n# Find the \"support\" bullet paragraph.\n$supportPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -like \"*The task will have support the following user inputs*\") {\n        $supportPara = $d.Paragraphs($i)\n        break\n    }\n}\n\n# Remove the existing (hidden) \"_GoBack\" bookmark wherever it currently is.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Re-add it inside the support bullet, splitting \"su\" | \"pport...\".\n$prefixLen = (\"The task will have su\").Length\n$newBookmarkPos = $supportPara.Range.Start + $prefixLen\n$newBookmarkRange = $d.Range($newBookmarkPos, $newBookmarkPos)\n$d.Bookmarks.Add(\"_GoBack\", $newBookmarkRange) | Out-Null\n"}
